$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the year table by one more column (2021 / 72) matching the
# formatting of the existing last column (Q).
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 72

# Update the view: select R1 (drops the stale topLeftCell scroll
# position and moves the active cell/selection to R1).
[void]$ws.Range("R1").Select()
